$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (advance by one day)
$ws.Range("A1").Value = 45311

# Update the price list values (fix for exceeded request bug / price scale correction)
$ws.Range("D28").Value = 534.582
$ws.Range("D29").Value = 552.874
$ws.Range("D30").Value = 562.708
$ws.Range("D31").Value = 579.591
$ws.Range("D32").Value = 830.009
$ws.Range("D33").Value = 800.441
$ws.Range("D34").Value = 1118.395
$ws.Range("D35").Value = 1147.937
